$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.392.47"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "2.285.22"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "503.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0956"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.91%  "
$ws.Range("E10").Value = "  +1.27%  "
$ws.Range("E11").Value = "  +3.45%  "
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "2.692.01"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.16%  "
$ws.Range("D15").Value = "54.340.64"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "2.290.23"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.84%  "
$ws.Range("E19").Value = "  +2.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "305.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.18%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "62.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "174.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.02%  "
$ws.Range("D30").Value = "0.0₃0693"
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.942"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.13%  "
$ws.Range("E36").Value = "  +1.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.04%  "
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("E40").Value = "  +2.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "124.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0497"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.53%  "
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.550"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "241.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.19%  "
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("E51").Value = "  +0.39%  "
